# Updates the cryptos price list (Sheet1) to the latest scraped snapshot:
# refreshed Price (D) / Volume(1h) (E) figures for most rows, plus a
# reordering of the FTXToken/HuobiToken rows (45/46) and the replacement
# of the FraxShare row (51) with THORChain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: Excel's COM Value setter auto-coerces plain numeric-looking
# strings (e.g. "248.92") into Double cells. The source data must stay
# text, so we briefly force a Text number format while assigning the
# value, then restore the default "Normal" style so no stray formatting
# is left behind on the cell.
function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) values
$ws.Cells.Item(2, 4).Value = "37.020.73"
$ws.Cells.Item(2, 5).Value = "  +1.45%  "
$ws.Cells.Item(3, 4).Value = "2.054.63"
$ws.Cells.Item(3, 5).Value = "  -2.25%  "
Set-TextValue 4 4 "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.19%  "
Set-TextValue 5 4 "248.92"
$ws.Cells.Item(5, 5).Value = "  -1.06%  "
Set-TextValue 6 4 "0.663"
$ws.Cells.Item(6, 5).Value = "  +0.81%  "
$ws.Cells.Item(7, 5).Value = "  -0.04%  "
Set-TextValue 8 4 "54.40"
$ws.Cells.Item(8, 5).Value = "  +14.24%  "
Set-TextValue 9 4 "60.77"
$ws.Cells.Item(9, 5).Value = "  +2.04%  "
$ws.Cells.Item(10, 5).Value = "  +1.03%  "
Set-TextValue 11 4 "0.0785"
$ws.Cells.Item(11, 5).Value = "  +5.49%  "
$ws.Cells.Item(12, 5).Value = "  +5.78%  "
Set-TextValue 13 4 "14.97"
$ws.Cells.Item(13, 5).Value = "  +3.43%  "
$ws.Cells.Item(14, 4).Value = "2.355.40"
$ws.Cells.Item(14, 5).Value = "  -2.28%  "
Set-TextValue 15 4 "0.811"
$ws.Cells.Item(15, 5).Value = "  -1.86%  "
Set-TextValue 16 4 "5.22"
$ws.Cells.Item(16, 5).Value = "  +2.51%  "
$ws.Cells.Item(17, 4).Value = "2.057.80"
$ws.Cells.Item(17, 5).Value = "  -2.05%  "
$ws.Cells.Item(18, 4).Value = "36.977.48"
$ws.Cells.Item(18, 5).Value = "  +1.35%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0942"
$ws.Cells.Item(19, 5).Value = "  +13.48%  "
Set-TextValue 20 4 "72.32"
$ws.Cells.Item(20, 5).Value = "  -0.73%  "
Set-TextValue 21 4 "14.19"
$ws.Cells.Item(21, 5).Value = "  +6.61%  "
$ws.Cells.Item(22, 5).Value = "  +3.23%  "
Set-TextValue 23 4 "236.28"
$ws.Cells.Item(23, 5).Value = "  -1.44%  "
Set-TextValue 24 4 "1.00"
$ws.Cells.Item(24, 5).Value = "  +0.01%  "
$ws.Cells.Item(25, 5).Value = "  -1.91%  "
Set-TextValue 26 4 "170.09"
$ws.Cells.Item(26, 5).Value = "  -0.62%  "
Set-TextValue 27 4 "8.99"
$ws.Cells.Item(27, 5).Value = "  -1.54%  "
Set-TextValue 28 4 "19.98"
$ws.Cells.Item(28, 5).Value = "  -6.83%  "
Set-TextValue 29 4 "1.98"
$ws.Cells.Item(29, 5).Value = "  -0.10%  "
$ws.Cells.Item(30, 5).Value = "  +0.34%  "
$ws.Cells.Item(31, 5).Value = "  +2.27%  "
$ws.Cells.Item(32, 5).Value = "  +0.63%  "
Set-TextValue 33 4 "1.04"
$ws.Cells.Item(33, 5).Value = "  +8.66%  "
Set-TextValue 34 4 "4.32"
$ws.Cells.Item(34, 5).Value = "  +6.30%  "
$ws.Cells.Item(35, 5).Value = "  -0.10%  "
Set-TextValue 36 4 "0.0858"
$ws.Cells.Item(36, 5).Value = "  -5.75%  "
$ws.Cells.Item(37, 5).Value = "  -2.79%  "
Set-TextValue 38 4 "1.76"
$ws.Cells.Item(38, 5).Value = "  -6.11%  "
$ws.Cells.Item(39, 5).Value = "  +0.11%  "
Set-TextValue 40 4 "0.105"
$ws.Cells.Item(40, 5).Value = "  +24.75%  "
Set-TextValue 41 4 "17.94"
$ws.Cells.Item(41, 5).Value = "  +11.66%  "
Set-TextValue 42 4 "0.0222"
$ws.Cells.Item(42, 5).Value = "  -0.05%  "
$ws.Cells.Item(43, 5).Value = "  -3.37%  "
Set-TextValue 44 4 "96.05"
$ws.Cells.Item(44, 5).Value = "  -1.58%  "
Set-TextValue 47 4 "2.39"
$ws.Cells.Item(47, 5).Value = "  +6.85%  "
$ws.Cells.Item(48, 4).Value = "1.291.63"
$ws.Cells.Item(48, 5).Value = "  -3.80%  "
$ws.Cells.Item(49, 5).Value = "  +2.96%  "
Set-TextValue 50 4 "12.96"
$ws.Cells.Item(50, 5).Value = "  -54.41%  "
# Row 45/46: HuobiToken/FTXToken swap positions with new values
$ws.Cells.Item(45, 2).Value = "FTXToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue 45 4 "4.25"
$ws.Cells.Item(45, 5).Value = "  +53.38%  "

$ws.Cells.Item(46, 2).Value = "HuobiToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue 46 4 "2.76"
$ws.Cells.Item(46, 5).Value = "  +0.72%  "

# Row 51: FraxShare replaced with THORChain
$ws.Cells.Item(51, 2).Value = "THORChain"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue 51 4 "4.05"
$ws.Cells.Item(51, 5).Value = "  +6.02%  "
